$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("D2").Value = 174
$ws.Range("D3").Value = 174
$ws.Range("D29").Value = 106

$ws.Range("D30").Select()
